$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "93.870.25"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.488.62"
$ws.Range("E3").Value = "  +5.04%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.65"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "624.83"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  +4.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.393"
$ws.Range("E8").Value = "  +3.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  +8.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.485.79"
$ws.Range("E11").Value = "  +4.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.06"
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("E13").Value = "  +5.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.26"
$ws.Range("E14").Value = "  +5.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.148.15"
$ws.Range("E15").Value = "  +5.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "93.682.26"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000250"
$ws.Range("E17").Value = "  +4.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.32"
$ws.Range("E18").Value = "  +5.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.496.66"
$ws.Range("E19").Value = "  +5.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.53"
$ws.Range("E20").Value = "  +15.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.89"
$ws.Range("E21").Value = "  +6.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.498"
$ws.Range("E22").Value = "  +12.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "518.90"
$ws.Range("E23").Value = "  +7.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.38"
$ws.Range("E24").Value = "  +5.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.69"
$ws.Range("E25").Value = "  +9.13%  "
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.96"
$ws.Range("E27").Value = "  +5.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.24"
$ws.Range("E28").Value = "  +6.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.668.97"
$ws.Range("E29").Value = "  +4.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.92"
$ws.Range("E30").Value = "  +11.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.42"
$ws.Range("E31").Value = "  +3.16%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  +3.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("E35").Value = "  +5.70%  "
$ws.Range("E36").Value = "  +5.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.558"
$ws.Range("E37").Value = "  +7.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "575.74"
$ws.Range("E38").Value = "  +10.79%  "
$ws.Range("E39").Value = "  +7.34%  "
$ws.Range("E40").Value = "  +4.05%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.922"
$ws.Range("E42").Value = "  +5.73%  "
$ws.Range("E43").Value = "  +2.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.75"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("E45").Value = "  +6.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.70"
$ws.Range("E46").Value = "  +2.98%  "
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.54"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.15"
$ws.Range("E49").Value = "  +3.36%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.21"
$ws.Range("E50").Value = "  +2.33%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.16"
$ws.Range("E51").Value = "  +2.94%  "
